$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "HK_R_acc_G"

$ws.Range("A2").Value = 86.227951153324284
$ws.Range("A3").Value = 84.396200814111268
$ws.Range("A4").Value = 84.871099050203526
$ws.Range("A5").Value = 83.175033921302571
$ws.Range("A6").Value = 87.449118046132966
$ws.Range("A7").Value = 89.620081411126179
$ws.Range("A8").Value = 87.788331071913163
$ws.Range("A9").Value = 87.856173677069208
$ws.Range("A10").Value = 87.924016282225239
$ws.Range("A11").Value = 87.924016282225239
$ws.Range("A12").Value = 89.552238805970148
$ws.Range("A13").Value = 89.687924016282224
$ws.Range("A14").Value = 88.059701492537314
$ws.Range("A15").Value = 88.12754409769336
$ws.Range("A16").Value = 88.059701492537314
$ws.Range("A17").Value = 87.652645861601087
$ws.Range("A18").Value = 89.891451831750331
$ws.Range("A19").Value = 89.891451831750331
$ws.Range("A20").Value = 85.413839891451843
$ws.Range("A21").Value = 80.936227951153327
$ws.Range("A22").Value = 85.61736770691995
$ws.Range("A23").Value = 84.73541383989145
$ws.Range("A24").Value = 85.278154681139753
$ws.Range("A25").Value = 84.803256445047495
$ws.Range("A26").Value = 87.584803256445042
$ws.Range("A27").Value = 87.652645861601087
$ws.Range("A28").Value = 87.788331071913163
$ws.Range("A29").Value = 88.534599728629587
$ws.Range("A30").Value = 88.602442333785618
$ws.Range("A31").Value = 88.805970149253739
$ws.Range("A32").Value = 81.682496607869737
$ws.Range("A33").Value = 85.549525101763919
$ws.Range("A34").Value = 85.753052917232026
$ws.Range("A35").Value = 83.175033921302571
$ws.Range("A36").Value = 83.175033921302571
$ws.Range("A37").Value = 86.36363636363636
$ws.Range("A38").Value = 83.039348710990495
$ws.Range("A39").Value = 83.10719131614654
$ws.Range("A40").Value = 85.006784260515602
$ws.Range("A41").Value = 87.584803256445042
$ws.Range("A42").Value = 87.042062415196739
$ws.Range("A43").Value = 87.24559023066486
$ws.Range("A44").Value = 86.770691994572587
$ws.Range("A45").Value = 86.499321573948436
$ws.Range("A46").Value = 86.36363636363636
$ws.Range("A47").Value = 86.635006784260511
$ws.Range("A48").Value = 83.310719131614647
$ws.Range("A49").Value = 86.770691994572587
